# Auto-generated Excel COM-interop script to update cryptos.xlsx data
# Mirrors the diff: refreshed Price (D) and Volume(1h) (E) columns,
# plus a reorder of the InternetComputer(DFINITY) / LidoDAOToken rows (28/29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
}

Set-TextValue ($ws.Cells.Item(2, 4)) '27.717.81'
Set-TextValue ($ws.Cells.Item(2, 5)) '  +0.49%  '
Set-TextValue ($ws.Cells.Item(3, 4)) '1.862.91'
Set-TextValue ($ws.Cells.Item(3, 5)) '  +0.60%  '
Set-TextValue ($ws.Cells.Item(4, 4)) '1.022'
Set-TextValue ($ws.Cells.Item(4, 5)) '  -1.16%  '
Set-TextValue ($ws.Cells.Item(5, 4)) '320.69'
Set-TextValue ($ws.Cells.Item(5, 5)) '  -0.36%  '
Set-TextValue ($ws.Cells.Item(6, 5)) '  -1.02%  '
Set-TextValue ($ws.Cells.Item(7, 4)) '0.4357'
Set-TextValue ($ws.Cells.Item(7, 5)) '  -0.95%  '
Set-TextValue ($ws.Cells.Item(8, 5)) '  +0.72%  '
Set-TextValue ($ws.Cells.Item(9, 4)) '0.07434'
Set-TextValue ($ws.Cells.Item(9, 5)) '  +0.24%  '
Set-TextValue ($ws.Cells.Item(10, 4)) '0.8839'
Set-TextValue ($ws.Cells.Item(10, 5)) '  +0.96%  '
Set-TextValue ($ws.Cells.Item(11, 4)) '21.61'
Set-TextValue ($ws.Cells.Item(11, 5)) '  +0.63%  '
Set-TextValue ($ws.Cells.Item(12, 4)) '1.866.22'
Set-TextValue ($ws.Cells.Item(12, 5)) '  +0.33%  '
Set-TextValue ($ws.Cells.Item(13, 4)) '6.749'
Set-TextValue ($ws.Cells.Item(13, 5)) '  +0.77%  '
Set-TextValue ($ws.Cells.Item(14, 4)) '5.493'
Set-TextValue ($ws.Cells.Item(14, 5)) '  -0.72%  '
Set-TextValue ($ws.Cells.Item(15, 4)) '0.07136'
Set-TextValue ($ws.Cells.Item(15, 5)) '  -1.15%  '
Set-TextValue ($ws.Cells.Item(16, 4)) '86.62'
Set-TextValue ($ws.Cells.Item(16, 5)) '  +4.40%  '
Set-TextValue ($ws.Cells.Item(17, 4)) '1.024'
Set-TextValue ($ws.Cells.Item(17, 5)) '  -1.17%  '
Set-TextValue ($ws.Cells.Item(18, 4)) '0.000009073'
Set-TextValue ($ws.Cells.Item(18, 5)) '  +0.50%  '
Set-TextValue ($ws.Cells.Item(19, 4)) '1.019'
Set-TextValue ($ws.Cells.Item(19, 5)) '  -1.07%  '
Set-TextValue ($ws.Cells.Item(20, 4)) '15.45'
Set-TextValue ($ws.Cells.Item(20, 5)) '  +0.27%  '
Set-TextValue ($ws.Cells.Item(21, 4)) '27.708.66'
Set-TextValue ($ws.Cells.Item(21, 5)) '  +0.41%  '
Set-TextValue ($ws.Cells.Item(22, 4)) '5.289'
Set-TextValue ($ws.Cells.Item(22, 5)) '  +0.56%  '
Set-TextValue ($ws.Cells.Item(23, 5)) '  -1.78%  '
Set-TextValue ($ws.Cells.Item(24, 4)) '2.090.06'
Set-TextValue ($ws.Cells.Item(24, 5)) '  +0.96%  '
Set-TextValue ($ws.Cells.Item(25, 4)) '2.039'
Set-TextValue ($ws.Cells.Item(25, 5)) '  +6.04%  '
Set-TextValue ($ws.Cells.Item(26, 4)) '157.50'
Set-TextValue ($ws.Cells.Item(26, 5)) '  -0.27%  '
Set-TextValue ($ws.Cells.Item(27, 4)) '18.73'
Set-TextValue ($ws.Cells.Item(27, 5)) '  +0.02%  '
Set-TextValue ($ws.Cells.Item(30, 4)) '120.55'
Set-TextValue ($ws.Cells.Item(30, 5)) '  +2.92%  '
Set-TextValue ($ws.Cells.Item(31, 4)) '0.09053'
Set-TextValue ($ws.Cells.Item(31, 5)) '  -0.14%  '
Set-TextValue ($ws.Cells.Item(32, 5)) '  +2.20%  '
Set-TextValue ($ws.Cells.Item(33, 4)) '0.7686'
Set-TextValue ($ws.Cells.Item(33, 5)) '  +0.85%  '
Set-TextValue ($ws.Cells.Item(34, 4)) '3.033'
Set-TextValue ($ws.Cells.Item(34, 5)) '  +5.13%  '
Set-TextValue ($ws.Cells.Item(35, 4)) '4.563'
Set-TextValue ($ws.Cells.Item(35, 5)) '  +1.26%  '
Set-TextValue ($ws.Cells.Item(36, 4)) '1.021'
Set-TextValue ($ws.Cells.Item(36, 5)) '  -1.02%  '
Set-TextValue ($ws.Cells.Item(37, 4)) '1.141'
Set-TextValue ($ws.Cells.Item(37, 5)) '  -0.58%  '
Set-TextValue ($ws.Cells.Item(38, 4)) '0.01980'
Set-TextValue ($ws.Cells.Item(38, 5)) '  +0.30%  '
Set-TextValue ($ws.Cells.Item(39, 4)) '0.05303'
Set-TextValue ($ws.Cells.Item(39, 5)) '  +0.07%  '
Set-TextValue ($ws.Cells.Item(40, 4)) '2.878'
Set-TextValue ($ws.Cells.Item(40, 5)) '  +2.53%  '
Set-TextValue ($ws.Cells.Item(41, 4)) '0.5190'
Set-TextValue ($ws.Cells.Item(41, 5)) '  +0.80%  '
Set-TextValue ($ws.Cells.Item(42, 4)) '6.946'
Set-TextValue ($ws.Cells.Item(42, 5)) '  +3.24%  '
Set-TextValue ($ws.Cells.Item(43, 5)) '  +0.31%  '
Set-TextValue ($ws.Cells.Item(44, 4)) '8.694'
Set-TextValue ($ws.Cells.Item(44, 5)) '  +2.59%  '
Set-TextValue ($ws.Cells.Item(45, 5)) '  +1.66%  '
Set-TextValue ($ws.Cells.Item(46, 4)) '110.12'
Set-TextValue ($ws.Cells.Item(46, 5)) '  +1.14%  '
Set-TextValue ($ws.Cells.Item(47, 4)) '1.715'
Set-TextValue ($ws.Cells.Item(47, 5)) '  +0.49%  '
Set-TextValue ($ws.Cells.Item(48, 4)) '1.021'
Set-TextValue ($ws.Cells.Item(48, 5)) '  -1.13%  '
Set-TextValue ($ws.Cells.Item(49, 4)) '0.06510'
Set-TextValue ($ws.Cells.Item(49, 5)) '  +1.82%  '
Set-TextValue ($ws.Cells.Item(50, 4)) '0.4719'
Set-TextValue ($ws.Cells.Item(50, 5)) '  +1.53%  '
Set-TextValue ($ws.Cells.Item(51, 4)) '1.870'
Set-TextValue ($ws.Cells.Item(51, 5)) '  +1.00%  '

# Rows 28/29: InternetComputer(DFINITY) and LidoDAOToken swap ranking order
Set-TextValue ($ws.Cells.Item(28, 2)) 'InternetComputer(DFINITY)'
Set-TextValue ($ws.Cells.Item(28, 3)) 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue ($ws.Cells.Item(28, 4)) '5.368'
Set-TextValue ($ws.Cells.Item(28, 5)) '  +2.04%  '

Set-TextValue ($ws.Cells.Item(29, 2)) 'LidoDAOToken'
Set-TextValue ($ws.Cells.Item(29, 3)) 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue ($ws.Cells.Item(29, 4)) '1.992'
Set-TextValue ($ws.Cells.Item(29, 5)) '  +1.10%  '
